# Insert a new weekly record above the last existing data row.
# The existing last row (50) holds the Sept-2021 reading which gets
# pushed down to row 51; row 50 becomes the new May-2022 reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the current last row (50) down to row 51 by inserting a new row
# at position 50 (Excel shifts row 50 -> 51, and everything below it).
# Insert() carries the date-format style on column D down through the
# shift, so both D50 (new) and D51 (old row 50's data) keep it.
$ws.Rows.Item(50).Insert()

# Fill the new row 50 with the new weekly record.
$ws.Cells.Item(50, 1).Value = 2
$ws.Cells.Item(50, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(50, 3).Value = "Coquimbo"
$ws.Cells.Item(50, 4).Value = 44706
$ws.Cells.Item(50, 5).Value = 4
$ws.Cells.Item(50, 6).Value = 100112022
$ws.Cells.Item(50, 7).Value = "Arveja Verde"
$ws.Cells.Item(50, 8).Value = "Perfection"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 160
$ws.Cells.Item(50, 11).Value = 25000
$ws.Cells.Item(50, 12).Value = 26000
$ws.Cells.Item(50, 13).Value = 25500
$ws.Cells.Item(50, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(50, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(50, 16).Value = 1020
$ws.Cells.Item(50, 17).Value = 25
$ws.Cells.Item(50, 18).Value = "Hortaliza"
